# Update Spain - Aragon "casos_coronavirus_hospitales" data.
# Appends one more day's worth of hospital records (2020-06-01, Excel
# serial date 43983) to the bottom of the data table, replicating the
# same 20-row hospital block structure used for the previous day
# (2020-05-31, rows 1153:1172, serial 43982), then fixes up the values
# that actually changed day over day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the previous day's 20-row block (A1153:H1172) into the new
# block (A1173:H1192) - first formats, then values, so cell styles and
# shared-string references both match the source block exactly.
$srcRange = $ws.Range("A1153:H1172")
$dstRange = $ws.Range("A1173:H1192")

$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$srcRange.Copy()
$dstRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# New block is the next day: 2020-06-01 (serial 43983).
$ws.Range("A1173:A1192").Value = 43983

# "camas_uci_ocupadas" (column C) values that changed from the prior day.
$ws.Range("C1173").Value = 8
$ws.Range("C1174").Value = 24
$ws.Range("C1177").Value = 2
$ws.Range("C1178").Value = 4
$ws.Range("C1179").Value = 3
$ws.Range("C1183").Value = 1
$ws.Range("C1184").Value = 8
$ws.Range("C1188").ClearContents()

$wb.Save()
